# Update the "Förändrad" (C column) date values from 2023-10-25 (45224)
# to 2023-11-03 (45233) for rows 2 through 10 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45233
    }
}
